# Update "想去人数" (interest count, column F) values across sheets to
# match the freshly generated gh-pages data output (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 71
$ws.Cells.Item(6, 6).Value = 610
$ws.Cells.Item(7, 6).Value = 325
$ws.Cells.Item(8, 6).Value = 527
$ws.Cells.Item(10, 6).Value = 10772
$ws.Cells.Item(14, 6).Value = 2040
$ws.Cells.Item(15, 6).Value = 890
$ws.Cells.Item(18, 6).Value = 209
$ws.Cells.Item(22, 6).Value = 127
$ws.Cells.Item(23, 6).Value = 207
$ws.Cells.Item(24, 6).Value = 688
$ws.Cells.Item(26, 6).Value = 221
$ws.Cells.Item(27, 6).Value = 2377
$ws.Cells.Item(28, 6).Value = 673
$ws.Cells.Item(29, 6).Value = 3165
$ws.Cells.Item(30, 6).Value = 1014
$ws.Cells.Item(31, 6).Value = 731
$ws.Cells.Item(33, 6).Value = 22
$ws.Cells.Item(39, 6).Value = 1
$ws.Cells.Item(41, 6).Value = 1263
$ws.Cells.Item(44, 6).Value = 130
$ws.Cells.Item(45, 6).Value = 223
$ws.Cells.Item(48, 6).Value = 4074
$ws.Cells.Item(49, 6).Value = 78

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 4094
$ws.Cells.Item(12, 6).Value = 356

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 411

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 411
$ws.Cells.Item(5, 6).Value = 71
$ws.Cells.Item(7, 6).Value = 610
$ws.Cells.Item(9, 6).Value = 10772
$ws.Cells.Item(14, 6).Value = 209
$ws.Cells.Item(17, 6).Value = 127
$ws.Cells.Item(18, 6).Value = 207
$ws.Cells.Item(19, 6).Value = 4094
$ws.Cells.Item(21, 6).Value = 688
$ws.Cells.Item(23, 6).Value = 221
$ws.Cells.Item(24, 6).Value = 673
$ws.Cells.Item(25, 6).Value = 3165
$ws.Cells.Item(26, 6).Value = 1014
$ws.Cells.Item(29, 6).Value = 731
$ws.Cells.Item(31, 6).Value = 22
$ws.Cells.Item(35, 6).Value = 1264
$ws.Cells.Item(38, 6).Value = 130
$ws.Cells.Item(39, 6).Value = 223
$ws.Cells.Item(44, 6).Value = 4074
$ws.Cells.Item(49, 6).Value = 78
